# Update the repeated block of rows 105-143 (columns A-J) with the new
# set of recalculated values, as produced by the model run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(3.279178071484681, 10.43797634328358, 19.67299221131186, 11.3450539176449, 15.94628430386411, 20.54751469008332, 25.14874508059701, 29.74997547111071, 34.3512058573299, 36.95856974328358)

for ($row = 105; $row -le 143; $row++) {
    for ($col = 1; $col -le 10; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
